# Update source/batch details to the new calibration cycle (Batch 42, Dec 10 2018)
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Reference Sheet Check")
$ws2 = $wb.Worksheets.Item("Check of Specific Date & time")

# --- "Reference Sheet Check" sheet -----------------------------------------

# Descriptive text block (A2:A5) - these are the "master" cells; the other
# sheet pulls them in via formulas, so updating them here is sufficient.
$ws1.Range("A2").Value = "Manufacurer's Calibration:383.8GBq @01:45 CET on 27 Nov 2018"
$ws1.Range("A3").Value = "Flexisource No.:  NLF 01 D85E-3112"
$ws1.Range("A4").Value = "CCSEO Batch Number: 42"
$ws1.Range("A5").Value = "Date of source installation @ CCSEO: Dec 10  2018"

# New manufacturer calibration date/time (A8) and installation date (E8)
$ws1.Range("A8").Value = 43431.072916666664
$ws1.Range("E8").Value = 43444.5

# New source activity at calibration (GBq)
$ws1.Range("A10").Value = 383.8

# --- "Check of Specific Date & time" sheet ----------------------------------

# Specific-date cell driving that sheet's decay calculation
$ws2.Range("E8").Value = 43466.5

# --- Restore the selected cell / active sheet -------------------------------

$ws1.Range("G10").Select()
$ws2.Activate()
